$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = -0.8456721019465476
$ws.Range("H2").Value = -0.001467681055260073
$ws.Range("B3").Value = -0.000009611129812764352
$ws.Range("H3").Value = -0.0009973959156024392
$ws.Range("B4").Value = -1.605341970103837
$ws.Range("H4").Value = -0.04899678965826126
$ws.Range("B5").Value = -0.001052723966949998
$ws.Range("D5").Value = -0.01489468452200526
$ws.Range("E5").Value = -0.003434566428040853
$ws.Range("H5").Value = -0.007501600139221409
$ws.Range("J5").Value = -0.00001156419847347934
$ws.Range("B6").Value = -1.452425980958651
$ws.Range("H6").Value = -0.002520714697269
$ws.Range("B7").Value = -0.02962907710752916
$ws.Range("H7").Value = -0.0009043117831879499
$ws.Range("B8").Value = -0.00001051225801518285
$ws.Range("C8").Value = -0.009929239795383182
$ws.Range("H8").Value = -0.01512774051298038
$ws.Range("B9").Value = -0.00007411219601749508
$ws.Range("H9").Value = 0.007358308512266376
$ws.Range("B10").Value = -0.000358142663799299
$ws.Range("D10").Value = 0.00101741758044227
$ws.Range("E10").Value = 0.0002346064170524187
$ws.Range("H10").Value = -0.0441112302005422
$ws.Range("J10").Value = 0.0000007899203939842891
$ws.Range("B11").Value = -0.0000295911003860283
$ws.Range("H11").Value = -0.1334753475530306
